$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.635.36'
$ws.Range('E2').Value = '  +1.41%  '

$ws.Range('D3').Value = '2.272.64'
$ws.Range('E3').Value = '  -1.89%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.74%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.20'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.99%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.569'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.84%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.19%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.86%  '

$ws.Range('E11').Value = '  -1.20%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.10'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.06%  '

$ws.Range('E13').Value = '  -1.50%  '

$ws.Range('D14').Value = '2.617.83'
$ws.Range('E14').Value = '  -1.95%  '

$ws.Range('D15').Value = '2.265.90'
$ws.Range('E15').Value = '  -1.67%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.26%  '

$ws.Range('D17').Value = '46.645.17'
$ws.Range('E17').Value = '  +1.69%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.796'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.97%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.50%  '

$ws.Range('D20').Value = '0.0₃0957'
$ws.Range('E20').Value = '  +1.63%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.20%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.74%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '248.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.94%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.29%  '

$ws.Range('E25').Value = '  +0.09%  '

$ws.Range('E26').Value = '  -5.57%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '41.59'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.58%  '

$ws.Range('E28').Value = '  -1.62%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.69'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.92%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.19'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.89%  '

$ws.Range('E31').Value = '  +7.58%  '

$ws.Range('B32').Value = 'LidoDAOToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +11.58%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '147.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.10%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0771'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.32%  '

$ws.Range('E36').Value = '  +7.65%  '

$ws.Range('E37').Value = '  -2.52%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.41%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.70'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.08%  '

$ws.Range('E40').Value = '  -5.14%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0296'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.14'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.23%  '

$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '91.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +15.15%  '

$ws.Range('D45').Value = '1.785.89'
$ws.Range('E45').Value = '  -0.36%  '

$ws.Range('E46').Value = '  -5.80%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '71.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.55%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.185'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.18%  '

$ws.Range('E49').Value = '  -0.48%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '94.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.95%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.20%  '
